# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for the rows
# whose automated dialog-act classification changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 14;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 20;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 39;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 44;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 47;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 48;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 63;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 65;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 69;  I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 75;  I = "ba"; J = "Appreciation" },
    @{ Row = 79;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 82;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 90;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 99;  I = "ba"; J = "Appreciation" },
    @{ Row = 100; I = "sv"; J = "Statement-opinion" },
    @{ Row = 104; I = "ba"; J = "Appreciation" },
    @{ Row = 112; I = "sv"; J = "Statement-opinion" },
    @{ Row = 113; I = "sv"; J = "Statement-opinion" },
    @{ Row = 117; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 122; I = "sv"; J = "Statement-opinion" },
    @{ Row = 124; I = "sv"; J = "Statement-opinion" },
    @{ Row = 129; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 141; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 146; I = "sv"; J = "Statement-opinion" },
    @{ Row = 147; I = "sv"; J = "Statement-opinion" },
    @{ Row = 149; I = "%";  J = "Uninterpretable" },
    @{ Row = 168; I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 173; I = "aa"; J = "Agree/Accept" },
    @{ Row = 174; I = "b";  J = "Acknowledge (Backchannel)" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
